$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '65.734.61'
$ws.Range('E2').Value = '  +0.93%  '

# Row 3
$ws.Range('D3').Value = '3.384.14'

# Row 4
$ws.Range('E4').Value = '  -0.07%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '563.58'
$ws.Range('E5').Value = '  +0.65%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '176.52'
$ws.Range('E6').Value = '  +0.55%  '

# Row 7
$ws.Range('E7').Value = '  +0.59%  '

# Row 8
$ws.Range('D8').Value = '3.377.85'
$ws.Range('E8').Value = '  -0.35%  '

# Row 9
$ws.Range('E9').Value = '  -0.02%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.174'
$ws.Range('E10').Value = '  +2.61%  '

# Row 11
$ws.Range('E11').Value = '  +0.44%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.87'
$ws.Range('E12').Value = '  -1.85%  '

# Row 14
$ws.Range('E14').Value = '  +1.20%  '

# Row 15
$ws.Range('D15').Value = '3.929.99'
$ws.Range('E15').Value = '  -0.46%  '

# Row 16
$ws.Range('B16').Value = 'TRON'
$ws.Range('C16').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.120'
$ws.Range('E16').Value = '  +0.38%  '

# Row 17
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.20'
$ws.Range('E17').Value = '  -1.03%  '

# Row 18
$ws.Range('D18').Value = '3.386.40'
$ws.Range('E18').Value = '  -0.81%  '

# Row 19
$ws.Range('D19').Value = '65.725.35'
$ws.Range('E19').Value = '  +0.94%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.89'
$ws.Range('E20').Value = '  +0.15%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.997'
$ws.Range('E21').Value = '  +0.25%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '462.18'
$ws.Range('E22').Value = '  -2.47%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.93'
$ws.Range('E23').Value = '  -0.87%  '

# Row 24
$ws.Range('E24').Value = '  +8.67%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '89.55'
$ws.Range('E25').Value = '  +2.68%  '

# Row 26
$ws.Range('E26').Value = '  -1.03%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.92'
$ws.Range('E27').Value = '  +0.16%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.64'
$ws.Range('E28').Value = '  -2.31%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.72'
$ws.Range('E29').Value = '  -1.05%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.08'
$ws.Range('E30').Value = '  -0.47%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.61'
$ws.Range('E31').Value = '  -1.56%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.48'
$ws.Range('E32').Value = '  -0.61%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '580.58'
$ws.Range('E33').Value = '  +1.25%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '62.27'
$ws.Range('E34').Value = '  +0.68%  '

# Row 35
$ws.Range('E35').Value = '  -0.35%  '

# Row 36
$ws.Range('E36').Value = '  +0.06%  '

# Row 37
$ws.Range('E37').Value = '  +1.76%  '

# Row 38
$ws.Range('E38').Value = '  +1.34%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.99'
$ws.Range('E39').Value = '  +0.44%  '

# Row 40
$ws.Range('E40').Value = '  +1.32%  '

# Row 41
$ws.Range('D41').Value = '0.0₃0746'
$ws.Range('E41').Value = '  -2.24%  '

# Row 42
$ws.Range('D42').Value = '3.107.33'
$ws.Range('E42').Value = '  +0.31%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.84'
$ws.Range('E43').Value = '  -1.06%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0417'
$ws.Range('E44').Value = '  +0.26%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.44'

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.16'
$ws.Range('E47').Value = '  -0.06%  '

# Row 48
$ws.Range('E48').Value = '  -0.05%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '140.96'
$ws.Range('E49').Value = '  +2.44%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.57'
$ws.Range('E50').Value = '  +8.87%  '

# Row 51
$ws.Range('B51').Value = 'WEMIXToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.56'
$ws.Range('E51').Value = '  -1.95%  '
